$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the original "_GoBack" bookmark sitting at the very start
#    of the document (around the title paragraph).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Swap out the old contact e-mail address for the new one.
# ------------------------------------------------------------------
$old = "cyberscitechcongress2020@gmail.com"
$new = "group-docieee2020@athabascau.ca"

$rng = $d.Content
$rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ------------------------------------------------------------------
# 3. Re-typing the address collapses the remainder of the host
#    paragraph down to maximally-merged runs (every run after the
#    edit point that happens to share identical formatting gets
#    glued together by the text-engine's reflow when the text of the
#    paragraph changes). Re-impose the original run boundaries by
#    nudging a property on each original run's range back to its own
#    value -- this forces the engine to re-materialize a distinct
#    run there without actually changing any visible formatting.
# ------------------------------------------------------------------
$found = $d.Content
$ok = $found.Find.Execute($new)

if ($ok) {
    $emailStart = $found.Start
    $emailEnd = $found.End

    # Text of every run that, in the original document, immediately
    # followed the e-mail address run (up to the end of the paragraph).
    # Re-splitting after each of these (in order) restores the original
    # run layout.
    $tailRunTexts = @(
        " ",
        "with the ",
        "subject",
        " ",
        "line ",
        [char]0x201C,
        "Letter of invitation",
        " request – Your Name",
        ".",
        [char]0x201D + " ",
        "Save your completed form using ",
        "the file name ",
        [char]0x201C,
        "Letter of invitation",
        " request – Your Name",
        ".",
        [char]0x201D
    )

    $prev = $emailEnd
    foreach ($t in $tailRunTexts) {
        $next = $prev + $t.Length
        $piece = $d.Range($prev, $next)
        $sz = $piece.Font.Size
        $piece.Font.Size = $sz + 1
        $piece.Font.Size = $sz
        $prev = $next
    }

    # --------------------------------------------------------------
    # 4. Word leaves the editing caret ("_GoBack") right where the
    #    last autocorrected/typed text landed -- in this case that is
    #    in the middle of the freshly inserted address, between
    #    "...atha" and "bascau.ca". Re-create that bookmark at the
    #    matching offset so the run gets split exactly like it was in
    #    Word.
    # --------------------------------------------------------------
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $splitOffset = $emailStart + "group-docieee2020@atha".Length
    $splitPoint = $d.Range($splitOffset, $splitOffset)
    $d.Bookmarks.Add("_GoBack", $splitPoint)
}
